$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns so numeric-looking
# strings (e.g. '1.011', '92.78', '0.00001110') are stored as text,
# matching the inlineStr cells in the source data (not coerced to numbers).
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '29.131.84'
$ws.Range('E2').Value = '  -1.05%  '
$ws.Range('D3').Value = '1.989.46'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('D4').Value = '1.011'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '329.87'
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('D6').Value = '1.009'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = '0.4981'
$ws.Range('E7').Value = '  -0.56%  '
$ws.Range('D8').Value = '0.4205'
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('D9').Value = '54.76'
$ws.Range('E9').Value = '  +3.87%  '
$ws.Range('D10').Value = '0.09287'
$ws.Range('E10').Value = '  +4.30%  '
$ws.Range('D11').Value = '1.098'
$ws.Range('E11').Value = '  -2.17%  '
$ws.Range('D12').Value = '22.99'
$ws.Range('E12').Value = '  -1.01%  '
$ws.Range('D13').Value = '2.001.50'
$ws.Range('E13').Value = '  -2.00%  '
$ws.Range('D14').Value = '7.979'
$ws.Range('E14').Value = '  -1.00%  '
$ws.Range('D15').Value = '6.463'
$ws.Range('E15').Value = '  -0.80%  '
$ws.Range('D16').Value = '1.011'
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value = '92.78'
$ws.Range('E17').Value = '  -3.37%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '0.00001110'
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('D19').Value = '0.06728'
$ws.Range('E19').Value = '  +1.55%  '
$ws.Range('D20').Value = '19.39'
$ws.Range('E20').Value = '  -1.49%  '
$ws.Range('D21').Value = '1.010'
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('D22').Value = '5.989'
$ws.Range('E22').Value = '  +0.38%  '
$ws.Range('D23').Value = '29.122.01'
$ws.Range('E23').Value = '  -1.21%  '
$ws.Range('D24').Value = '12.03'
$ws.Range('E24').Value = '  +1.30%  '
$ws.Range('D25').Value = '2.284'
$ws.Range('E25').Value = '  +1.06%  '
$ws.Range('D26').Value = '2.265.68'
$ws.Range('E26').Value = '  +0.99%  '
$ws.Range('D27').Value = '20.85'
$ws.Range('E27').Value = '  +1.19%  '
$ws.Range('D28').Value = '156.52'
$ws.Range('E28').Value = '  -1.50%  '
$ws.Range('D29').Value = '6.336'
$ws.Range('E29').Value = '  -2.11%  '
$ws.Range('D30').Value = '2.267'
$ws.Range('E30').Value = '  -2.88%  '
$ws.Range('D31').Value = '127.78'
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('D32').Value = '1.054'
$ws.Range('E32').Value = '  +0.80%  '
$ws.Range('D33').Value = '0.09848'
$ws.Range('E33').Value = '  -0.87%  '
$ws.Range('D34').Value = '1.543'
$ws.Range('E34').Value = '  -1.46%  '
$ws.Range('D35').Value = '5.828'
$ws.Range('D36').Value = '3.741'
$ws.Range('E36').Value = '  -1.21%  '
$ws.Range('D37').Value = '0.02429'
$ws.Range('E37').Value = '  -1.20%  '
$ws.Range('D38').Value = '1.322'
$ws.Range('E38').Value = '  +2.76%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '9.151'
$ws.Range('E39').Value = '  -4.33%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '0.06404'
$ws.Range('E40').Value = '  +1.01%  '
$ws.Range('D41').Value = '0.6504'
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').Value = '11.52'
$ws.Range('E42').Value = '  -1.58%  '
$ws.Range('D43').Value = '0.2005'
$ws.Range('E43').Value = '  -2.95%  '
$ws.Range('D44').Value = '1.010'
$ws.Range('E44').Value = '  +0.23%  '
$ws.Range('D45').Value = '0.6233'
$ws.Range('E45').Value = '  -1.57%  '
$ws.Range('E46').Value = '  +4.97%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '2.194'
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '13.41'
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').Value = '3.488'
$ws.Range('E49').Value = '  -0.95%  '
$ws.Range('D50').Value = '0.00000000331'
$ws.Range('E50').Value = '  +1.83%  '
$ws.Range('D51').Value = '0.07008'
$ws.Range('E51').Value = '  +0.10%  '
